$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.35844298829246668,
    -0.0099999997172801613,
    -0.0089999997127687692,
    -0.011999999923464344,
    -0.0059999997173534325,
    -0.005999999709231929,
    -0.019999999655000877,
    -0.019999999654384482,
    -0.0059999997084414503,
    -0.0059999997091111368,
    -0.0044999997146710768,
    -0.0059999997092843316,
    -0.0059999997119195569,
    -0.011999999690046614,
    -0.0059999997148008077,
    0.04472543497619208,
    -0.0059999997177850872,
    0.013265226269796315,
    -0.0089999997203653592,
    -0.0089999997128415998,
    -0.008999999711184703,
    -0.055590928186678212,
    -0.0089999997102774287,
    -0.041999999578637315,
    -0.041999999576296965,
    -0.0059999997082620382,
    -0.0059999997065243171,
    -0.0059999996998394423,
    -0.011999999672175576,
    -0.019999999639154886,
    -0.014999999655213259,
    -0.020999999631610144,
    -0.0059999996891288987
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
